$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.692.64"
$ws.Range("E2").Value = "  -0.80%  "
$ws.Range("D3").Value = "2.570.68"
$ws.Range("E3").Value = "  +0.38%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "581.30"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.49%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "143.21"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -3.03%  "
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("E8").Value = "  +0.45%  "
$ws.Range("E9").Value = "  -2.97%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "5.58"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.53%  "
$ws.Range("E11").Value = "  -0.43%  "
$ws.Range("E12").Value = "  -2.28%  "
$ws.Range("E13").Value = "  -2.69%  "
$ws.Range("D14").Value = "3.033.03"
$ws.Range("E14").Value = "  +0.56%  "
$ws.Range("D15").Value = "62.667.77"
$ws.Range("E15").Value = "  -0.67%  "
$ws.Range("E16").Value = "  -3.27%  "
$ws.Range("D17").Value = "2.580.66"
$ws.Range("E17").Value = "  +0.80%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "11.04"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -3.05%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "340.33"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.40%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.32"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.29%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.62"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.77%  "
$ws.Range("E22").Value = "  -0.03%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "66.84"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.33%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.57"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -4.31%  "
$ws.Range("E25").Value = "  -3.80%  "
$ws.Range("E26").Value = "  +0.84%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.999"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.10%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.88"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -2.35%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.20"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -4.24%  "
$ws.Range("E30").Value = "  -4.65%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "453.54"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +2.52%  "
$ws.Range("E32").Value = "  -3.95%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "176.77"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.05%  "
$ws.Range("E34").Value = "  +0.35%  "
$ws.Range("E35").Value = "  +0.10%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.399"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.96%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "18.82"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.57%  "
$ws.Range("E38").Value = "  -1.57%  "
$ws.Range("E39").Value = "  -0.01%  "
$ws.Range("E40").Value = "  -3.31%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "39.92"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.51%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "156.78"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +3.81%  "
$ws.Range("E43").Value = "  -4.04%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.632"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +3.12%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "20.99"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.67%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0532"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -3.99%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0958"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.81%  "
$ws.Range("E48").Value = "  -3.19%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "17.92"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.67%  "
$ws.Range("E50").Value = "  +0.26%  "
$ws.Range("E51").Value = "  -4.62%  "
